$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("D2").Value = "model.pkl"
$ws.Range("D3").Value = "model.pkl"
$ws.Range("D4").Value = "model.pkl"
$ws.Range("D5").Value = "model.xlsx"

$ws.Range("D6").Select()
